$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.228.35"
$ws.Range("E2").Value = "  -2.80%  "
$ws.Range("D3").Value = "1.549.35"
$ws.Range("E3").Value = "  -4.81%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'206.89"
$ws.Range("E5").Value = "  -3.40%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.476"
$ws.Range("E7").Value = "  -5.55%  "
$ws.Range("D8").Value = "'0.0607"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "1.764.76"
$ws.Range("E12").Value = "  -4.80%  "
$ws.Range("D13").Value = "1.545.16"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("E14").Value = "  -4.82%  "
$ws.Range("D15").Value = "'0.503"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "25.196.61"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "0.0₃0706"
$ws.Range("E17").Value = "  -4.42%  "
$ws.Range("D18").Value = "'58.50"
$ws.Range("E18").Value = "  -4.58%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "'185.65"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").Value = "'9.24"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("E23").Value = "  -4.15%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "'0.128"
$ws.Range("E25").Value = "  -4.14%  "
$ws.Range("D26").Value = "'139.56"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("E27").Value = "  -5.09%  "
$ws.Range("D28").Value = "'14.79"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("E29").Value = "  -5.19%  "
$ws.Range("E30").Value = "  -6.67%  "
$ws.Range("E31").Value = "  -4.69%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("E33").Value = "  -4.82%  "
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("D36").Value = "1.084.82"
$ws.Range("E36").Value = "  -3.62%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "'0.494"
$ws.Range("E39").Value = "  -5.38%  "
$ws.Range("D40").Value = "'2.24"
$ws.Range("E40").Value = "  -7.70%  "
$ws.Range("D41").Value = "'0.760"
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("D43").Value = "'92.55"
$ws.Range("D44").Value = "'5.03"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "1.680.07"
$ws.Range("E45").Value = "  -4.79%  "
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  +13.82%  "
$ws.Range("D47").Value = "'1.45"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").Value = "'52.22"
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("E49").Value = "  -5.93%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -2.06%  "
